# Etapas.xlsx: remove the stray numeric values that had been entered in
# column C (rows 2-5) underneath the "cordinador" header - the column is
# otherwise unused or populated elsewhere, so the cells are cleared
# entirely rather than just blanked out.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C5").ClearContents()
